$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.732.08"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "2.216.64"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.18"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.56"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.49"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.94"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.88"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "2.549.64"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.67"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "2.208.32"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.802"
$ws.Range("E18").Value = "  -4.10%  "
$ws.Range("D19").Value = "42.561.77"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.87"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.09"
$ws.Range("E23").Value = "  -9.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "228.94"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.14"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.38"
$ws.Range("E28").Value = "  -6.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.81"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.46"
$ws.Range("E32").Value = "  +13.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.26"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0792"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.64"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.50"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.79"
$ws.Range("E43").Value = "  -5.36%  "
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0984"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.70"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.422"
$ws.Range("E51").Value = "  +14.34%  "
